# Replicate the previous "HateCrimes.xls" merge step, this time for
# Offense = 'Hate Crime - Non Forcible Sex Offenses', sum2013 only.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sectors = @(
    "Public, 4-year or above",
    "Private nonprofit, 4-year or above",
    "Private for-profit, 4-year or above",
    "Public, 2-year",
    "Private nonprofit, 2-year",
    "Private for-profit, 2-year",
    "Public, less-than 2-year",
    "Private nonprofit, less-than 2-year",
    "Private for-profit, less-than 2-year"
)

$locations = @(
    "On Campus (excluding Residence Halls)",
    "On Campus (Residence Halls)",
    "Non-Campus",
    "Public Property"
)

$offense = "Hate Crime - Non Forcible Sex Offenses"
$date = "sum2013"

$row = 2081
foreach ($loc in $locations) {
    foreach ($sec in $sectors) {
        $ws.Cells.Item($row, 1).Value = $sec
        $ws.Cells.Item($row, 2).Value = $loc
        $ws.Cells.Item($row, 3).Value = $offense
        $ws.Cells.Item($row, 4).Value = $date
        $ws.Cells.Item($row, 5).Value = 0
        $row = $row + 1
    }
}

# Restore the view state (scroll position / active cell) to match the
# state left behind after entering this block of data.
$ws.Range("A2066").Select()
$ws.Range("F2084").Select()
